$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# Revision date and revision number bump (row 5)
$ws.Range("E5").Value = 46001
$ws.Range("G5").Value = 1.3

# Assign Risk Priority Number (RPN) values to risks in column G (rows 9-23)
$ws.Range("G9").Value = 10
$ws.Range("G10").Value = 2
$ws.Range("G11").Value = 11
$ws.Range("G12").Value = 3
$ws.Range("G13").Value = 4
$ws.Range("G14").Value = 12
$ws.Range("G15").Value = 1
$ws.Range("G16").Value = 5
$ws.Range("G17").Value = 6
$ws.Range("G18").Value = 13
$ws.Range("G19").Value = 15
$ws.Range("G20").Value = 7
$ws.Range("G21").Value = 8
$ws.Range("G22").Value = 14
$ws.Range("G23").Value = 9

# Vertically center the Description column text for risk rows
$ws.Range("C9:C23").VerticalAlignment = -4108

# Rows 20-23 used a left-aligned / locked variant of the description style;
# align it with the rest of the block (general horizontal alignment, unlocked)
$ws.Range("C20:C23").HorizontalAlignment = 1
$ws.Range("C20:C23").Locked = $false

# Move selection to reflect where the author was last working
$ws.Range("F14").Select() | Out-Null
